$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 currently holds the text "R40"; update it to the text "1".
$ws.Range("B11").Value = "1"
